$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-02-17 04:52:19"
$wsZhCn.Range("G5").Value = "2016-02-17 04:53:02"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-02-17 04:52:29"
$wsDeDe.Range("G5").Value = "2016-02-17 04:53:20"
